# Update the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) values on the zh-cn and de-de
# report sheets, as a newer handback report run produced later timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-26 09:35:44"
$wsZhCn.Range("G2").Value = "2016-01-26 09:36:27"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-26 09:35:59"
$wsDeDe.Range("G2").Value = "2016-01-26 09:36:48"
